# Regenerate merged AHB files
#  - rename the _old/_new header suffixes to _FV2210/_FV2304
#  - add a freeze pane under the header row
#  - turn the used range into an Excel table (Table1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row (row 1): "_old" -> "_FV2210", "_new" -> "_FV2304"
[void]$ws.Cells.Replace("_old", "_FV2210", 2)
[void]$ws.Cells.Replace("_new", "_FV2304", 2)

# 2) Freeze the header row (split under row 1)
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3) Convert the used range A1:U88 into a table named Table1
$range = $ws.Range("A1:U88")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"
